$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.04756558637480458
$ws.Range("C2").Value = 0.9986011772539274
$ws.Range("D2").Value = 0.1602021758772121
$ws.Range("G2").Value = 0.3552643978832445
$ws.Range("H2").Value = 0.9990000000000001

$ws.Range("B3").Value = 0.07600487841318371
$ws.Range("C3").Value = 0.9992790596668165
$ws.Range("D3").Value = 0.205298418927665
$ws.Range("G3").Value = 0.3552643978832445
$ws.Range("H3").Value = 0.9990000000000001

$ws.Range("B4").Value = 0.08843381759302323
$ws.Range("C4").Value = 0.9988026606982543
$ws.Range("D4").Value = 0.2397858522798944
$ws.Range("G4").Value = 0.3552643978832445
$ws.Range("H4").Value = 0.9990000000000001
